$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Range("A1").Interior.TintAndShade = 0.4
$newSheet.Range("A1").Interior.ThemeColor = 5
